$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column C: shift the week's dates forward by two weeks ---
$ws.Range("C2").Value = "07/22/2019"
$ws.Range("C3").Value = "07/23/2019"
$ws.Range("C4").Value = "07/24/2019"
$ws.Range("C5").Value = "07/25/2019"
$ws.Range("C6").Value = "07/26/2019"

# --- Column D: replace the daily status remarks with the new cycle-2 entries ---
$ws.Range("D2").Value = "1.`tReviewed Bug Snag errors for both PMall mobile and desktop site for dev environment. Created 3 new defects and assigned to Vinod for review.`n2.`tVerified the ticket #7950 and shared a minor observation with you on email. Kindly review and suggest.`n3.`tObserved large number of errors appearing in browser console for both desktop and mobile site. Shared screenshots in another email. Please review and suggest.`n4.`tAttended daily status meeting to discuss project updates, tasks , observations and requirements."

$ws.Range("D3").Value = "1.`tVerified the new functionalities on dev environment: #7799: Website Logo Update(desktop and mobile site) and #7873:Desktop Product Category Featured Redirected Store/Category Items. Shared observations with you and dev team for review.`n2.`tCreated/Updated test cases for #7799 and #7873.`n3.`tCreated status report for June release-cycle 2 which includes most of the tickets tested previously that are to be regressed or retested because of new changes, along with some new tickets deployed in July. Please find attached updated status report with Ticket-Wise Test Summary- Cycle 2, Test Case Execution Summary- Cycle 2, Test cases and Defects created in June release till date. "

$ws.Range("D4").Value = "1.`tRetested defects assigned to me in PMall admin: #7890, #7803, #7817, #7652 and #6990.`n2.`tPerformed further testing on #7799: Website logo update and #7539: Brand logo & Artist Page. Multiple issues are appearing on the branding page as discussed in call. I will discontinue testing on the same until further intimation.`n3.`tAttended daily standup. Discussed observations raised on 23rd July along with other issues. Mailed all issues to Sonny and Vinod on separate emails. Will verify the same tomorrow and raise the defects in PMall admin if any of them persist. Also, had discussion with Vinod for the scope of #7873.`n4.`tUpdated status report for June release-cycle 2. Please find updated sheet attached."

$ws.Range("D5").Value = "1.`tRetested  defect #7734 on dev environment assigned to me in PMall admin. Working fine on dev environment.`n2.`tCompleted task # 4410: Implement Facebook DPA Tags in PMall admin.`n3.`tPerformed testing on Gift Checkout Options on both mobile and desktop site on chrome browser. Shared observation documents with you and other developers on a separate email.`n4.`tCreated 7 new defects from #7980 to #7986.`n5.`tUpdated status report for June release-cycle 2. Please find updated sheet attached."

$ws.Range("D6").Value = "1.`tPerformed functional testing for remaining tickets today on Chrome Browser on Dev environment. Testing is pending for all tickets on different browsers and devices. Testing for #7539: Brand logo & Artist Page is on hold for now as per yesterday’s discussion.`n•`t#7304: Mobile Checkout Progress Bar`n•`t#7571: Product Page Thumbnail Image Enhancement`n•`t#7645: Mobile strike through fix`n•`t#6380 : #Add 301 Redirects`n•`t#7383: Wedding Page Refresh with new Template`n•`t#6791: Mobile round corner block refresh`n2.`tReviewed BugSnag errors for both mobile and desktop.`n3.`tRetested defects #7908, #7843 and #7891 on dev environment.`n4.`tAttended daily status meeting to discuss all observations found today and yesterday. Created 19 new defects from #7994 – #8012 including the bug snag errors and other console errors discussed today.`n5.`tUpdated status report for June release-cycle 2. Please find updated sheet attached."

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 120
$ws.Rows.Item(3).RowHeight = 150
$ws.Rows.Item(4).RowHeight = 150
$ws.Rows.Item(5).RowHeight = 120
$ws.Rows.Item(6).RowHeight = 240
